# Add a new column P "Estimated Start Date" to Sheet1 (dimension grows from O123 to P123).
# Header gets the same bold/centered/bordered style as the other header cells (copied from A1).
# Rows with a known estimated start date get a numeric Excel date serial, formatted like the
# existing date columns (style copied from M2, "YYYY-MM-DD HH:MM:SS"). Remaining rows get a blank
# placeholder (no estimated date yet), matching the source data export.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header ---
$ws.Range("P1").Value = "Estimated Start Date"
$ws.Range("A1").Copy()
$ws.Range("P1").PasteSpecial(-4122)   # xlPasteFormats: reuse the existing header style, no new style created
$excel.CutCopyMode = $false

# --- Rows with an estimated start date ---
$ws.Range("P4").Value = 45614.3841370088
$ws.Range("P8").Value = 45626.3841370088
$ws.Range("P25").Value = 45627.3841370088
$ws.Range("P26").Value = 45625.3841370088
$ws.Range("P29").Value = 45625.3841370088
$ws.Range("P33").Value = 45621.3841370088
$ws.Range("P38").Value = 45626.3841370088
$ws.Range("P54").Value = 45620.3841370088
$ws.Range("P55").Value = 45618.3841370088
$ws.Range("P63").Value = 45617.3841370088
$ws.Range("P64").Value = 45613.3841370088
$ws.Range("P68").Value = 45620.3841370088
$ws.Range("P92").Value = 45615.3841370088
$ws.Range("P93").Value = 45617.3841370088
$ws.Range("P94").Value = 45625.3841370088
$ws.Range("P95").Value = 45621.3841370088
$ws.Range("P100").Value = 45615.3841370088
$ws.Range("P108").Value = 45626.3841370088
$ws.Range("P115").Value = 45624.3841370088
$ws.Range("P119").Value = 45613.3841370088
$ws.Range("P121").Value = 45615.3841370088
$ws.Range("P122").Value = 45625.3841370088
$ws.Range("P123").Value = 45618.3841370088

$dateRows = @(4,8,25,26,29,33,38,54,55,63,64,68,92,93,94,95,100,108,115,119,121,122,123)
$ws.Range("M2").Copy()
foreach ($r in $dateRows) {
    $ws.Range("P$r").PasteSpecial(-4122)   # xlPasteFormats: reuse the existing date style, no new style created
}
$excel.CutCopyMode = $false

# --- Remaining rows: no estimated start date yet ---
$emptyRows = @(2,3,5,6,7,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,27,28,30,31,32,34,35,36,37,39,40,41,42,43,44,45,46,47,48,49,50,51,52,53,56,57,58,59,60,61,62,65,66,67,69,70,71,72,73,74,75,76,77,78,79,80,81,82,83,84,85,86,87,88,89,90,91,96,97,98,99,101,102,103,104,105,106,107,109,110,111,112,113,114,116,117,118,120)
foreach ($r in $emptyRows) {
    $ws.Range("P$r").Value = ""
}
